$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" "261.54"
Set-TextValue $ws "E2" "0.27%"
Set-TextValue $ws "D3" "26.65"
Set-TextValue $ws "E3" "-2.80%"
Set-TextValue $ws "D4" "4.714"
Set-TextValue $ws "E4" "0.52%"
Set-TextValue $ws "D5" "0.06172"
Set-TextValue $ws "E5" "1.30%"
Set-TextValue $ws "D6" "6.707"
Set-TextValue $ws "E6" "0.65%"
Set-TextValue $ws "D7" "0.8504"
Set-TextValue $ws "E7" "0.17%"
Set-TextValue $ws "D8" "0.9123"
Set-TextValue $ws "E8" "-1.73%"
Set-TextValue $ws "D9" "0.1408"
Set-TextValue $ws "E9" "0.23%"
Set-TextValue $ws "D10" "0.05160"
Set-TextValue $ws "E10" "5.56%"
Set-TextValue $ws "D11" "0.07098"
Set-TextValue $ws "E11" "-0.08%"
Set-TextValue $ws "D12" "0.03118"
Set-TextValue $ws "E12" "1.31%"
Set-TextValue $ws "D13" "0.09042"
Set-TextValue $ws "E13" "-0.21%"
Set-TextValue $ws "D14" "0.001539"
Set-TextValue $ws "E14" "-0.14%"
Set-TextValue $ws "D15" "0.0006167"
Set-TextValue $ws "E15" "1.40%"
Set-TextValue $ws "D16" "0.006079"
Set-TextValue $ws "E16" "-0.67%"
Set-TextValue $ws "D17" "3.451"
Set-TextValue $ws "E17" "0.10%"
Set-TextValue $ws "D18" "3.172"
Set-TextValue $ws "E18" "0.76%"
Set-TextValue $ws "D19" "2.187"
Set-TextValue $ws "E19" "1.11%"
Set-TextValue $ws "D21" "0.1301"
Set-TextValue $ws "E21" "-0.28%"
Set-TextValue $ws "D22" "4.104"
Set-TextValue $ws "E22" "0.58%"
Set-TextValue $ws "D23" "0.04232"
Set-TextValue $ws "E23" "-0.15%"
Set-TextValue $ws "D24" "0.001176"
Set-TextValue $ws "E24" "-3.75%"
Set-TextValue $ws "E25" "6.55%"
Set-TextValue $ws "E26" "0.05%"
Set-TextValue $ws "E27" "4.11%"
Set-TextValue $ws "D40" "0.03990"
Set-TextValue $ws "E40" "3.48%"
Set-TextValue $ws "D41" "0.1111"
Set-TextValue $ws "E41" "-0.24%"
Set-TextValue $ws "D42" "0.004145"
Set-TextValue $ws "E42" "1.50%"
Set-TextValue $ws "E43" "-3.33%"
Set-TextValue $ws "D44" "0.01326"
Set-TextValue $ws "E44" "-18.28%"
Set-TextValue $ws "E45" "0.46%"
Set-TextValue $ws "E46" "0.05%"
Set-TextValue $ws "E47" "-61.04%"
Set-TextValue $ws "D48" "0.2571"
Set-TextValue $ws "E48" "89.72%"
Set-TextValue $ws "D49" "0.00002102"
Set-TextValue $ws "E49" "0.05%"
Set-TextValue $ws "E50" "0.05%"
